# BASIC_STAT hp, attackDamage, magicDamage -> long long
# Update the type row (row 2) for HP / AttackDamage / MagicDamage columns
# on both the Item_BasicAddStat and Item_BasicMulStat sheets from "int" to
# "long long", and mirror the author's resulting selection / active-sheet
# state.

$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("Item_BasicAddStat")
$wsMul = $wb.Worksheets.Item("Item_BasicMulStat")

# --- data changes: HP (C), AttackDamage (E), MagicDamage (F) -> "long long"
$wsAdd.Range("C2").Value = "long long"
$wsAdd.Range("E2").Value = "long long"
$wsAdd.Range("F2").Value = "long long"

$wsMul.Range("C2").Value = "long long"
$wsMul.Range("E2").Value = "long long"
$wsMul.Range("F2").Value = "long long"

# --- selection / active sheet changes left behind by the editing session
$wsAdd.Activate() | Out-Null
$wsAdd.Range("A2:K2").Select() | Out-Null

$wsMul.Activate() | Out-Null
$wsMul.Range("A2:K2").Select() | Out-Null
